# Auto-applies the cell-level value changes described by the target diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 809
$ws.Range("F5").Value = 138
$ws.Range("F7").Value = 716
$ws.Range("F8").Value = 716
$ws.Range("F11").Value = 21
$ws.Range("F12").Value = 1099
$ws.Range("F13").Value = 854
$ws.Range("F14").Value = 700
$ws.Range("F17").Value = 1335
$ws.Range("F20").Value = 531
$ws.Range("F21").Value = 553
$ws.Range("F25").Value = 1059
$ws.Range("F29").Value = 455
$ws.Range("F34").Value = 259
$ws.Range("F37").Value = 1229
$ws.Range("F40").Value = 3853

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("F25").Value = 232
$ws.Range("F29").Value = 228
$ws.Range("F33").Value = 1702

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("F8").Value = 956

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("F7").Value = 956
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "2024-03-09"
$ws.Range("B8").ClearFormats()
$ws.Range("C8").Value = "上海·S·CGE动漫游戏嘉年华"
$ws.Range("D8").Value = "军工路1076号 纪希片场(秀场)"
$ws.Range("E8").Value = "2024.03.09 10:00-03.10 17:00"
$ws.Range("F8").Value = 6744
$ws.Range("G8").Value = 80
$ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=81173"
$ws.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202401/TYA5FLkE1705891815532.jpeg"
$ws.Range("C9").Value = "上海·爱乐之城音乐会"
$ws.Range("D9").Value = "南京西路1376号 上海商城剧院"
$ws.Range("E9").Value = "2024.03.09 14:00-03.09 15:30"
$ws.Range("F9").Value = 37
$ws.Range("G9").Value = 168
$ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=81289"
$ws.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202401/ZZXtDrwZ1705996679699.jpeg"
$ws.Range("C10").Value = "上海·第五十三届燃梦星辰国潮动漫嘉年华-随机宅舞"
$ws.Range("D10").Value = "周家嘴路3608号 宝龙旭辉广场"
$ws.Range("E10").Value = "2024.03.09 10:20-03.10 16:30"
$ws.Range("F10").Value = 809
$ws.Range("G10").Value = 58
$ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=80571"
$ws.Range("I10").Value = "//i0.hdslb.com/bfs/openplatform/202401/SHH70VXN1704700240858.jpeg"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "2024-03-10"
$ws.Range("B11").ClearFormats()
$ws.Range("C11").Value = "上海·三森铃子10周年纪念2024演唱会"
$ws.Range("D11").Value = "宜昌路179号 万代南梦宫上海文化中心"
$ws.Range("E11").Value = "2024.03.10 18:00-03.10 19:30"
$ws.Range("F11").Value = 748
$ws.Range("G11").Value = 399
$ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=81433"
$ws.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202401/L8rmm2h81706236781799.jpeg"
$ws.Range("C12").Value = "上海·次元裂缝-X 新春二次元DJ派对"
$ws.Range("D12").Value = "海潮路133号B1 JUMP工坊"
$ws.Range("E12").Value = "2024.03.10 14:00-03.10 19:00"
$ws.Range("F12").Value = 138
$ws.Range("G12").Value = 60
$ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=81959"
$ws.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202402/MaO7WWLr1708482746780.jpeg"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "2024-03-15"
$ws.Range("B13").ClearFormats()
$ws.Range("C13").Value = "上海·坏孩纸物语の第35届动漫节之全民宅舞"
$ws.Range("D13").Value = "泸定路388号 桃源π商业广场"
$ws.Range("E13").Value = "2024.03.15 11:30-03.16 16:00"
$ws.Range("F13").Value = 4
$ws.Range("G13").Value = 20
$ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=82477"
$ws.Range("I13").Value = "//i2.hdslb.com/bfs/openplatform/202403/WHufQNn91709782559844.png"
$ws.Range("F14").Value = 716
$ws.Range("F15").Value = 716
$ws.Range("F17").Value = 21
$ws.Range("F18").Value = 1099
$ws.Range("F19").Value = 854
$ws.Range("F21").Value = 700
$ws.Range("F26").Value = 1335
$ws.Range("F29").Value = 531
$ws.Range("F30").Value = 553
$ws.Range("C32").Value = "上海·第七届ACBC动漫盛典-平金中心国漫浪潮嘉年华"
$ws.Range("D32").Value = "七莘路1599弄(七莘路地铁站1号口步行60米) 平金中心"
$ws.Range("E32").Value = "2024.03.30 10:00-03.31 18:00"
$ws.Range("F32").Value = 3
$ws.Range("G32").Value = 48.8
$ws.Range("H32").Value = "https://show.bilibili.com/platform/detail.html?id=82487"
$ws.Range("I32").Value = "//i0.hdslb.com/bfs/openplatform/202403/yszE1z2O1709794608587.png"
$ws.Range("F34").Value = 1059
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "2024-04-04"
$ws.Range("B36").ClearFormats()
$ws.Range("C36").Value = "上海·原神X星穹铁道ONLY"
$ws.Range("D36").Value = "逸仙路301号靠纪念路路口 上海宝丰联大酒店"
$ws.Range("E36").Value = "2024.04.04 10:00-04.04 17:00"
$ws.Range("F36").Value = 723
$ws.Range("G36").Value = 60
$ws.Range("H36").Value = "https://show.bilibili.com/platform/detail.html?id=80299"
$ws.Range("I36").Value = "//i2.hdslb.com/bfs/openplatform/202312/V0xu26Cl1703753850690.jpeg"
$ws.Range("C37").Value = "上海·首届sunshine跨次元动漫游戏展"
$ws.Range("D37").Value = "莘福路288号 美莘商业广场"
$ws.Range("E37").Value = "2024.04.04 10:00-04.05 17:00"
$ws.Range("F37").Value = 521
$ws.Range("G37").Value = 50
$ws.Range("H37").Value = "https://show.bilibili.com/platform/detail.html?id=82417"
$ws.Range("I37").Value = "//i0.hdslb.com/bfs/openplatform/202403/y4v1H69x1709708980441.jpeg"
$ws.Range("C38").Value = "上海·魔都coser动漫展-C展"
$ws.Range("D38").Value = "海潮路133号B1 JUMP工坊"
$ws.Range("F38").Value = 455
$ws.Range("G38").Value = 60
$ws.Range("H38").Value = "https://show.bilibili.com/platform/detail.html?id=82104"
$ws.Range("I38").Value = "//i2.hdslb.com/bfs/openplatform/202402/vqcswCGV1708942084553.jpeg"
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "2024-04-05"
$ws.Range("B39").ClearFormats()
$ws.Range("C39").Value = "上海·怀旧番ONLY"
$ws.Range("D39").Value = "逸仙路270号  上海宝丰联大酒店"
$ws.Range("E39").Value = "2024.04.05 10:00-04.05 17:00"
$ws.Range("F39").Value = 449
$ws.Range("H39").Value = "https://show.bilibili.com/platform/detail.html?id=80575"
$ws.Range("I39").Value = "//i1.hdslb.com/bfs/openplatform/202401/y4uWdyPT1704700763902.jpeg"
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "2024-04-06"
$ws.Range("B40").ClearFormats()
$ws.Range("C40").Value = "上海·从Butter-Fly到夏目之爱してる —— “好想大声说爱你”动漫钢琴演奏会"
$ws.Range("D40").Value = "复兴中路1380号 捷豹上海交响音乐厅"
$ws.Range("E40").Value = "2024.04.06 19:30-04.06 21:30"
$ws.Range("F40").Value = 39
$ws.Range("G40").Value = 80
$ws.Range("H40").Value = "https://show.bilibili.com/platform/detail.html?id=80050"
$ws.Range("I40").Value = "//i0.hdslb.com/bfs/openplatform/202312/0iJP3TY61703056498448.jpeg"
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "2024-04-13"
$ws.Range("B41").ClearFormats()
$ws.Range("C41").Value = "上海·《四月是你的谎言》——“公生”与“薰”的钢琴小提琴唯美经典音乐集"
$ws.Range("D41").Value = "丁香路425号 上海东方艺术中心"
$ws.Range("E41").Value = "2024.04.13 19:30-04.13 21:30"
$ws.Range("F41").Value = 247
$ws.Range("H41").Value = "https://show.bilibili.com/platform/detail.html?id=78667"
$ws.Range("I41").Value = "//i1.hdslb.com/bfs/openplatform/202311/bTP7w6GD1700130122940.jpeg"
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "2024-04-13"
$ws.Range("B42").ClearFormats()
$ws.Range("C42").Value = "上海·第三届奇卡波利国潮嘉年华"
$ws.Range("D42").Value = "漕宝路3366号 七宝万科广场"
$ws.Range("E42").Value = "2024.04.13 10:30-04.14 16:30"
$ws.Range("F42").Value = 7
$ws.Range("G42").Value = 28.8
$ws.Range("H42").Value = "https://show.bilibili.com/platform/detail.html?id=82376"
$ws.Range("I42").Value = "//i1.hdslb.com/bfs/openplatform/202403/64i0bjSy1709692398951.jpeg"
$ws.Range("F44").Value = 259
$ws.Range("F46").Value = 228
$ws.Range("F47").Value = 1702
$ws.Range("F48").Value = 1702
$ws.Range("F49").Value = 1229
$ws.Range("F51").Value = 3853
